$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (C) column was refreshed for every remaining data row (2-27):
# 45265 -> 45266 (2023-12-05 -> 2023-12-06)
$ws.Range("C2:C27").Value = 45266

# Row 28 (A 61490-2023) was removed from the log entirely
$ws.Rows(28).Delete()

# Row 27 no longer carries an explicit custom row height - let it revert
# to the sheet's default height
$ws.Rows(27).AutoFit()
